$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- ReceiptItems table header: "For Whome" -> "Ownership" ---
$ws1.Range("D14").Value = "Ownership"

# --- Receipts table header row 8: shift "FK - GroupID" from H8 to new I8,
#     rename H8 to "Scanned", add a blank K8 (style match, extends table) ---
$ws1.Range("H8").Copy()
$ws1.Range("I8").PasteSpecial(-4122)
$ws1.Range("I8").Value = "FK - GroupID"
$ws1.Range("H8").Value = "Scanned"

$ws1.Range("J8").Copy()
$ws1.Range("K8").PasteSpecial(-4122)

# --- Row 9 (receipt 1): old H9 numeric value (1) moves to I9;
#     H9 becomes a boolean "Scanned" flag (FALSE) ---
$ws1.Range("H9").Copy()
$ws1.Range("I9").PasteSpecial(-4122)
$ws1.Range("I9").Value = 1
$ws1.Range("H9").Value = $false

$ws1.Range("J9").Copy()
$ws1.Range("K9").PasteSpecial(-4122)

# --- Row 10 (receipt 2): new I10 numeric value (2);
#     H10 becomes a boolean "Scanned" flag (TRUE) ---
$ws1.Range("H10").Copy()
$ws1.Range("I10").PasteSpecial(-4122)
$ws1.Range("I10").Value = 2
$ws1.Range("H10").Value = $true

$ws1.Range("J10").Copy()
$ws1.Range("K10").PasteSpecial(-4122)

# --- Selection moves to G21 (scrolled down to the scanned-receipt workings area) ---
$ws1.Range("G21").Select()
